# Update "想去人数" (interest count) figures for several 漫展 events
# across the "展览" and "全部类型" sheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 244
$wsExhibit.Range("F6").Value = 10018
$wsExhibit.Range("F10").Value = 5530
$wsExhibit.Range("F17").Value = 297
$wsExhibit.Range("F22").Value = 1513

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 244
$wsAll.Range("F7").Value = 10018
$wsAll.Range("F11").Value = 5530
$wsAll.Range("F18").Value = 297
$wsAll.Range("F23").Value = 1513
